# Commit: feat: add 2022-Q1 data
#
# The workbook's last sheet ("总计", a running summary of all quarters) is
# repurposed to hold the new "2022-Q1" per-fund detail table (matching the
# layout used by every other quarterly sheet), and a brand new "总计" sheet
# is appended at the end with the summary table refreshed to include the
# new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Turn the current last sheet ("总计") into the new "2022-Q1" detail
#    sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Name = "2022-Q1"

# -- header row ---------------------------------------------------------
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"

# new header cells E1:H1 need the same style as the existing header cells
$q1.Range("B1").Copy()
$q1.Range("E1").PasteSpecial(-4122)
$q1.Range("E1").Value = "股票总仓位"

$q1.Range("B1").Copy()
$q1.Range("F1").PasteSpecial(-4122)
$q1.Range("F1").Value = "仓位占比"

$q1.Range("B1").Copy()
$q1.Range("G1").PasteSpecial(-4122)
$q1.Range("G1").Value = "持有市值(亿元)"

$q1.Range("B1").Copy()
$q1.Range("H1").PasteSpecial(-4122)
$q1.Range("H1").Value = "仓位排名"

# -- fund rows ------------------------------------------------------------
# columns D,E,F,G hold numeric-looking figures that are stored as *text*
# (matching every other quarterly sheet), column H is a real number.
$funds = @(
    @(0,  "513090", "易方达中证香港证券投资主题ETF",               "11.07", "96.47", "17.35", "1.9206", 1),
    @(1,  "501050", "华夏沪港通上证50AH优选指数（LOF）A",           "25.94", "92.28", "2.88",  "0.7471", 7),
    @(2,  "004496", "前海开源多元策略灵活配置混合A",                "3.09",  "91.25", "7.91",  "0.2444", 5),
    @(3,  "004497", "前海开源多元策略灵活配置混合C",                "2.06",  "91.25", "7.91",  "0.1629", 5),
    @(4,  "011722", "前海开源深圳特区精选股票型证券投资基金A",        "3.00",  "85.40", "4.53",  "0.1359", 10),
    @(5,  "011355", "华泰柏瑞港股通时代机遇混合型证券投资基金A",      "1.13",  "90.93", "6.46",  "0.0730", 5),
    @(6,  "003413", "华泰柏瑞新经济沪港深灵活配置混合",              "0.54",  "92.57", "5.92",  "0.0320", 6),
    @(7,  "011356", "华泰柏瑞港股通时代机遇混合型证券投资基金C",      "0.40",  "90.93", "6.46",  "0.0258", 5),
    @(8,  "011723", "前海开源深圳特区精选股票型证券投资基金C",        "0.44",  "85.40", "4.53",  "0.0199", 10),
    @(9,  "501067", "招商富时中国A-H50指数（LOF）A",               "0.21",  "94.63", "3.42",  "0.0072", 7),
    @(10, "006395", "华夏沪港通上证50AH优选指数（LOF）C",           "0.25",  "92.28", "2.88",  "0.0072", 7),
    @(11, "501068", "招商富时中国A-H50指数（LOF）C",               "0.05",  "94.63", "3.42",  "0.0017", 7),
    @(12, "160922", "大成恒生综合中小型股指数(QDII-LOF)A",          "0.10",  "92.44", "1.03",  "0.0010", 10),
    @(13, "002860", "前海开源沪港深新机遇灵活配置混合",              "0.01",  "83.26", "6.61",  "0.0007", 6)
)

$rowNum = 2
foreach ($fund in $funds) {
    if ($rowNum -gt 6) {
        # rows 7..15 are brand new -- clone column A's style from row 2
        $q1.Range("A2").Copy()
        $q1.Cells.Item($rowNum, 1).PasteSpecial(-4122)
    }
    $q1.Cells.Item($rowNum, 1).Value = $fund[0]

    # fund code must stay text (leading zeros, e.g. "004496")
    $q1.Cells.Item($rowNum, 2).NumberFormat = "@"
    $q1.Cells.Item($rowNum, 2).Value = $fund[1]

    $q1.Cells.Item($rowNum, 3).Value = $fund[2]

    $q1.Cells.Item($rowNum, 4).NumberFormat = "@"
    $q1.Cells.Item($rowNum, 4).Value = $fund[3]

    $q1.Cells.Item($rowNum, 5).NumberFormat = "@"
    $q1.Cells.Item($rowNum, 5).Value = $fund[4]

    $q1.Cells.Item($rowNum, 6).NumberFormat = "@"
    $q1.Cells.Item($rowNum, 6).Value = $fund[5]

    $q1.Cells.Item($rowNum, 7).NumberFormat = "@"
    $q1.Cells.Item($rowNum, 7).Value = $fund[6]

    $q1.Cells.Item($rowNum, 8).Value = $fund[7]

    $rowNum = $rowNum + 1
}

# ---------------------------------------------------------------------
# 2. Append a brand new "总计" sheet after "2022-Q1" with the refreshed
#    summary table (adds the 2022-Q1 row on top, shifting the rest down).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# copy the header style from "2022-Q1" onto the new header cells
$q1.Range("B1").Copy()
$total.Range("B1").PasteSpecial(-4122)
$q1.Range("B1").Copy()
$total.Range("C1").PasteSpecial(-4122)
$q1.Range("B1").Copy()
$total.Range("D1").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$summary = @(
    @(0, "2022-Q1", 14, 3.38),
    @(1, "2021-Q4", 13, 3.64),
    @(2, "2021-Q3", 20, 3.5),
    @(3, "2021-Q2", 14, 2.88),
    @(4, "2021-Q1", 12, 2.33),
    @(5, "2020-Q4", 21, 5.26)
)

$rowNum = 2
foreach ($row in $summary) {
    $q1.Range("A2").Copy()
    $total.Cells.Item($rowNum, 1).PasteSpecial(-4122)
    $total.Cells.Item($rowNum, 1).Value = $row[0]
    $total.Cells.Item($rowNum, 2).Value = $row[1]
    $total.Cells.Item($rowNum, 3).Value = $row[2]
    $total.Cells.Item($rowNum, 4).Value = $row[3]
    $rowNum = $rowNum + 1
}
